$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read existing values (B1:E4) into a local array before we overwrite anything
$data = $ws.Range("B1:E4").Value2

# Write them back shifted one column to the left (A1:D4)
$ws.Range("A1:D4").Value2 = $data

# Clear the now-stale column E which used to hold data
$ws.Range("E1:E4").ClearContents()

$ws.Range("A1:A1048576").Select()
